$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the BTec logo picture (held in both page headers) from
# "image1.jpg" to "image2.jpg". Selecting the shape first and then
# renaming it through $word.Selection avoids a stale-handle issue that
# hits InlineShape.Name when the shape lives in a header/footer range.
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                [void]$shp.Select()
                $word.Selection.InlineShapes.Item(1).Name = "image2.jpg"
            }
        }
    }
}

# Rename the Pearson Edexcel logo picture (held in both page footers)
# from "image2.png" to "image1.png".
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                [void]$shp.Select()
                $word.Selection.InlineShapes.Item(1).Name = "image1.png"
            }
        }
    }
}
